# Insert a new data row at row 537 (pushing the existing rows 537:562 down
# to 538:563) and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 537..562 down to 538..563, copying formatting (date style on
# column D, etc.) from the row above - this matches Excel's default
# "insert row" behaviour.
$ws.Rows.Item(537).Insert()

# Populate the newly inserted row 537 with the new record.
$ws.Cells.Item(537, 1).Value = 7
$ws.Cells.Item(537, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(537, 3).Value = "Ñuble"
$ws.Cells.Item(537, 4).Value = 45041
$ws.Cells.Item(537, 5).Value = 16
$ws.Cells.Item(537, 6).Value = 100114001
$ws.Cells.Item(537, 7).Value = "Papa"
$ws.Cells.Item(537, 8).Value = "Patagonia"
$ws.Cells.Item(537, 9).Value = "1a (cosecha)"
$ws.Cells.Item(537, 10).Value = 500
$ws.Cells.Item(537, 11).Value = 10000
$ws.Cells.Item(537, 12).Value = 11000
$ws.Cells.Item(537, 13).Value = 10500
$ws.Cells.Item(537, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(537, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(537, 16).Value = 420
$ws.Cells.Item(537, 17).Value = 25
$ws.Cells.Item(537, 18).Value = "Hortaliza"
